$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.142.33'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.29%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.210.55'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.91%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.44%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.07%  '

$ws.Range('E7').Value = '  +0.10%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.208.07'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.94%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.531'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.97%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.156'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.89%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.55'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.02%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.483'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.41%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000260'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.91%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.55'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.75%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.726.51'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.27%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.105.23'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.37%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.202.20'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.20%  '

$ws.Range('E18').Value = '  +0.63%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.12'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.20%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '488.61'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.70%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.41%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.726'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.47%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.36%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.29%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.49%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.995'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.59%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.19%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.84'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.12%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.132'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +38.81%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.31'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.15%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.87'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.05%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.78'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.83%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.26'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.70%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.29%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.12'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.42%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.20%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.96'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.10%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.30'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.70%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '479.97'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.60%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0748'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.94%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0410'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.70%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.127'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.27%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.58'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.97%  '

$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.49'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.17%  '

$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.951.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.04%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.284'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.00%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '27.84'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.33%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.39'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.84%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.998'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.05%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.117'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.13%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '121.07'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.01%  '
